# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F3").Value = 300
$wsExpo.Range("F4").Value = 2786
$wsExpo.Range("F6").Value = 594

# Sheet 4: 全部类型 (All types)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value = 300
$wsAll.Range("F6").Value = 2786
$wsAll.Range("F8").Value = 594
